{"js": "// Update the three \"Lista de Necessidades\" requirement bullets (N01-N03)\n// and force the section's page orientation to explicit Portrait.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map the old bullet text (as a prefix) to its full replacement text.\nconst replacements = [\n  {\n    prefix: \"N01:\",\n    text: \"N01: Website para divulga\u00e7\u00e3o dos contatos e servi\u00e7os prestados.\"\n  },\n  {\n    prefix: \"N02:\",\n    text: \"N02: Gest\u00e3o de estoque das pe\u00e7as dispon\u00edveis.\"\n  },\n  {\n    prefix: \"N03:\",\n    text: \"N03: Gest\u00e3o de OSs.\"\n  }\n];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const current = para.text || \"\";\n  const match = replacements.find((r) => current.indexOf(r.prefix) === 0);\n  if (match) {\n    // \"Replace\" keeps the paragraph/run formatting and just swaps the text.\n    para.insertText(match.text, \"Replace\");\n  }\n}\nawait context.sync();\n\n// Explicitly set the page orientation to Portrait so <w:pgSz> gets an\n// explicit w:orient=\"portrait\" attribute.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < sections.items.length; i++) {\n  sections.items[i].pageSetup.orientation = \"Portrait\";\n}\nawait context.sync();\n", "ps1": "# Update the three \"Lista de Necessidades\" requirement bullets (N01-N03)\n# and force the section's page orientation to explicit Portrait.\n\n$d = $word.ActiveDocument\n\n# Map the old bullet prefix to the new full paragraph text.\n$replacements = @{\n    \"N01:\" = \"N01: Website para divulga\u00e7\u00e3o dos contatos e servi\u00e7os prestados.\"\n    \"N02:\" = \"N02: Gest\u00e3o de estoque das pe\u00e7as dispon\u00edveis.\"\n    \"N03:\" = \"N03: Gest\u00e3o de OSs.\"\n}\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    foreach ($prefix in $replacements.Keys) {\n        if ($t.StartsWith($prefix)) {\n            # Assigning .Text keeps the paragraph/run formatting and just\n            # swaps the text content (paragraph mark is untouched).\n            $p.Range.Text = $replacements[$prefix]\n            break\n        }\n    }\n}\n\n# Explicitly set the page orientation to Portrait (wdOrientPortrait = 0) so\n# <w:pgSz> gets an explicit w:orient=\"portrait\" attribute.\n$d.PageSetup.Orientation = 0\n"}
